# Generate Report for Handoff
# Updates the Status and timestamp cells on the Overview / zh-cn / de-de
# sheets to reflect a fresh handoff, and shrinks the now-narrower Status /
# language columns to fit the shorter "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "Handed back: in sync with en-US" -> "Ready for handoff"
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value = "Ready for handoff"
$dede.Range("C2").Value = "Ready for handoff"

# --- Timestamps for the new handoff
$overview.Range("G2").Value = "2016-10-19 11:41:52"
$dede.Range("H2").Value = "2016-10-19 11:41:52"
$zhcn.Range("H2").Value = "2016-10-19 11:41:41"

# --- Narrow the Status / language columns to fit the shorter text
# (engine snaps ColumnWidth to a pixel grid; 16.33 is the input that lands
# closest to the recorded 17.2159881591797 post-snap)
$overview.Range("E:E").ColumnWidth = 16.33
$overview.Range("F:F").ColumnWidth = 16.33
$zhcn.Range("C:C").ColumnWidth = 16.33
$dede.Range("C:C").ColumnWidth = 16.33
